$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-12-22 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-23 Monday", 2)
$d.Content.Find.Execute("692÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "942÷3=", 2)
$d.Content.Find.Execute("713÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "563÷4=", 2)
$d.Content.Find.Execute("831÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "896÷5=", 2)
$d.Content.Find.Execute("447÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "226÷3=", 2)
$d.Content.Find.Execute("312÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "393÷4=", 2)
$d.Content.Find.Execute("273÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "267÷5=", 2)
$d.Content.Find.Execute("619÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "822÷4=", 2)
$d.Content.Find.Execute("780÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "975÷7=", 2)
$d.Content.Find.Execute("528÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "868÷4=", 2)
$d.Content.Find.Execute("954÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "550÷9=", 2)
$d.Content.Find.Execute("987÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "316÷7=", 2)
$d.Content.Find.Execute("365÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "842÷3=", 2)
$d.Content.Find.Execute("785÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "874÷5=", 2)
$d.Content.Find.Execute("265÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "132÷5=", 2)
$d.Content.Find.Execute("904÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "619÷8=", 2)
$d.Content.Find.Execute("341÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "820÷8=", 2)
$d.Content.Find.Execute("541÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "930÷3=", 2)
$d.Content.Find.Execute("249÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "531÷4=", 2)
$d.Content.Find.Execute("319÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "137÷4=", 2)
$d.Content.Find.Execute("176÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "263÷2=", 2)
$d.Content.Find.Execute("108÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "306÷2=", 2)
$d.Content.Find.Execute("732÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "985÷4=", 2)
$d.Content.Find.Execute("903÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "642÷2=", 2)
$d.Content.Find.Execute("159÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "401÷9=", 2)
$d.Content.Find.Execute("759÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "602÷3=", 2)
